{"js": "// Remove the stray \"tirm\" run at the very start of the document,\n// and collapse each \"<id>\" / \"p126v_N\" / \"</id>\" run triple into a\n// single run containing the combined text \"<id>p126v_N</id>\".\n\nconst body = context.document.body;\n\n// 1) Delete the leftover \"tirm\" text (first run of the document).\nconst tirmResults = body.search(\"tirm\", { matchCase: true });\ntirmResults.load(\"text\");\nawait context.sync();\nif (tirmResults.items.length > 0) {\n  tirmResults.items[0].delete();\n}\n\n// 2) Merge the split \"<id>...</id>\" runs back into a single run for\n// each occurrence found in the document.\nconst ids = [\"p126v_1\", \"p126v_2\", \"p126v_3\"];\nfor (let i = 0; i < ids.length; i++) {\n  const needle = \"<id>\" + ids[i] + \"</id>\";\n  const found = body.search(needle, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(needle, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the stray \"tirm\" run at the very start of the document, and\n# collapse each \"<id>\" / \"p126v_N\" / \"</id>\" run triple into a single\n# run containing the combined text \"<id>p126v_N</id>\".\n\n$d = $word.ActiveDocument\n\n# 1) Delete the leftover \"tirm\" text (first run of the document).\n$rngTirm = $d.Content\nif ($rngTirm.Find.Execute(\"tirm\")) {\n    $rngTirm.Delete()\n}\n\n# 2) Merge the split \"<id>...</id>\" runs back into a single run for\n# each of the three occurrences in the document (\"p126v_1\"..\"p126v_3\"),\n# searching forward each time so every occurrence is handled in turn.\n$searchFrom = 0\nforeach ($n in 1..3) {\n    $idValue = \"p126v_$n\"\n\n    $rng = $d.Range($searchFrom, $d.Content.End)\n    $found = $rng.Find.Execute(\"<id>\")\n    if (-not $found) { continue }\n\n    # $rng now covers exactly the \"<id>\" run. Compute the absolute\n    # boundaries of the following \"p126v_N\" and \"</id>\" runs.\n    $midStart = $rng.End\n    $midEnd = $midStart + $idValue.Length\n    $closeEnd = $midEnd + 5   # length of \"</id>\"\n\n    # Remove the \"p126v_N\" and \"</id>\" runs entirely.\n    $rngRest = $d.Range($midStart, $closeEnd)\n    $rngRest.Delete()\n\n    # Re-append the removed text onto the end of the \"<id>\" run so it\n    # merges back into a single run with that run's formatting.\n    $rng.Collapse(0)\n    $rng.InsertAfter($idValue + \"</id>\")\n\n    $searchFrom = $rng.End\n}\n"}
